$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1:H1 with same style as existing headers (copy from E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy formatting from existing header (A1:E1) to the new headers
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill F2:H20 with boolean FALSE values
$boolRange = $ws.Range("F2:H20")
$boolRange.Value = $false

$wb.Save()
